$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 1 cells with new values
$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -2,3)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: 4,-10)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: 8,8)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -5,2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: 4,6)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -6,6)"

# Add new column G with a new entry
$ws.Range("G1").Value = "(308051846, Eyal  Sofer: -3,-3)"

# Update cost and time rows
$ws.Range("A3").Value = "cost: 670.2245000844722"
$ws.Range("A4").Value = "time: 80.65306251055902"
